$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the data table with 5 new daily rows (101-105), matching the date-style
# formatting (style index carrying the custom date number format) used by column A.
$ws.Range("A100").Copy() | Out-Null
$ws.Range("A101:A105").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 101
$ws.Cells.Item(101, 1).Value2 = 45531
$ws.Cells.Item(101, 2).Value2 = 604.4175061499999
$ws.Cells.Item(101, 3).Value2 = 170.3822929455
$ws.Cells.Item(101, 9).Value2 = 262.819142649
$ws.Cells.Item(101, 11).Value2 = 285.447828102555
$ws.Cells.Item(101, 14).Value2 = 32.23453685088
$ws.Cells.Item(101, 15).Value2 = 0.021395996
$ws.Cells.Item(101, 17).Value2 = [double]"1.8336e-06"
$ws.Cells.Item(101, 21).Value2 = 254.3228880466148
$ws.Cells.Item(101, 26).Value2 = 184.24006325417

# Row 102
$ws.Cells.Item(102, 1).Value2 = 45532
$ws.Cells.Item(102, 2).Value2 = 600.550821069
$ws.Cells.Item(102, 3).Value2 = 175.3051737955
$ws.Cells.Item(102, 9).Value2 = 257.734284245
$ws.Cells.Item(102, 11).Value2 = 286.171545656385
$ws.Cells.Item(102, 14).Value2 = 31.29159806976
$ws.Cells.Item(102, 15).Value2 = 0.02154444
$ws.Cells.Item(102, 17).Value2 = [double]"1.8456e-06"
$ws.Cells.Item(102, 21).Value2 = 247.0309340130851
$ws.Cells.Item(102, 26).Value2 = 183.948774616614

# Row 103
$ws.Cells.Item(103, 1).Value2 = 45533
$ws.Cells.Item(103, 2).Value2 = 603.8479305180999
$ws.Cells.Item(103, 3).Value2 = 175.2552516235
$ws.Cells.Item(103, 9).Value2 = 250.662034
$ws.Cells.Item(103, 11).Value2 = 289.064606666355
$ws.Cells.Item(103, 14).Value2 = 32.19961319232
$ws.Cells.Item(103, 15).Value2 = 0.021520368
$ws.Cells.Item(103, 17).Value2 = [double]"1.836e-06"
$ws.Cells.Item(103, 21).Value2 = 247.2867920493493
$ws.Cells.Item(103, 26).Value2 = 190.211480324068

# Row 104
$ws.Cells.Item(104, 1).Value2 = 45534
$ws.Cells.Item(104, 2).Value2 = 601.4571167119
$ws.Cells.Item(104, 3).Value2 = 175.1436201
$ws.Cells.Item(104, 9).Value2 = 247.063243369
$ws.Cells.Item(104, 11).Value2 = 289.245998712
$ws.Cells.Item(104, 14).Value2 = 33.30552904672
$ws.Cells.Item(104, 15).Value2 = 0.021472224
$ws.Cells.Item(104, 17).Value2 = [double]"1.8648e-06"
$ws.Cells.Item(104, 21).Value2 = 248.3102241944061
$ws.Cells.Item(104, 26).Value2 = 192.957916049596

# Row 105
$ws.Cells.Item(105, 1).Value2 = 45535
$ws.Cells.Item(105, 2).Value2 = 599.9311952118999
$ws.Cells.Item(105, 3).Value2 = 174.2429409135
$ws.Cells.Item(105, 9).Value2 = 242.336473585
$ws.Cells.Item(105, 11).Value2 = 285.088692338865
$ws.Cells.Item(105, 14).Value2 = 31.30323928928
$ws.Cells.Item(105, 15).Value2 = 0.021379948
$ws.Cells.Item(105, 17).Value2 = [double]"1.8216e-06"
$ws.Cells.Item(105, 21).Value2 = 243.3209924872542
$ws.Cells.Item(105, 26).Value2 = 189.087938436352
